# Regenerate the handback-status report for the new handoff/handback run.
# Old source files:
#   2a460be5-ad8a-4081-84b0-fdb30031e67f.md
#   8804289f-0354-4f2b-ad11-01460922f2e2.md
# New source files:
#   224a5908-ba5e-498b-a323-6945e0c4307d.md
#   ffff16dbc55d-a242-4824-9f55-a0c20595bee7.md

$wb = $excel.ActiveWorkbook

$oldFile1 = "2a460be5-ad8a-4081-84b0-fdb30031e67f.md"
$oldFile2 = "8804289f-0354-4f2b-ad11-01460922f2e2.md"
$newFile1 = "224a5908-ba5e-498b-a323-6945e0c4307d.md"
$newFile2 = "ffff16dbc55d-a242-4824-9f55-a0c20595bee7.md"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("A3").Value = $newFile2

$wsOverview.Range("G2").Value = "2016-09-03 01:07:11"
$wsOverview.Range("G3").Value = "2016-09-03 01:07:11"

# Rebuild the two hyperlinks on column B with the new display text while
# keeping the same external targets (only the file name segment changes).
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0362aca070505afce837ff1ed03df425643cd9b7/e2e/$newFile1", "", "", "e2e\$newFile1")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0362aca070505afce837ff1ed03df425643cd9b7/e2e/$newFile2", "", "", "e2e\$newFile2")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("I3").Value = $newFile2

$zhXlf = "224a5908-ba5e-498b-a323-6945e0c4307d.974b29d69231ec61ba9d9827bd1e821688ad796b.zh-cn.xlf"
$wsZhCn.Range("G2").Value = $zhXlf
$wsZhCn.Range("J2").Value = $zhXlf
$wsZhCn.Range("G3").Value = $zhXlf
$wsZhCn.Range("J3").Value = $zhXlf

$wsZhCn.Range("H2").Value = "2016-09-03 01:07:01"
$wsZhCn.Range("H3").Value = "2016-09-03 01:07:01"
$wsZhCn.Range("K2").Value = "2016-09-03 01:07:30"
$wsZhCn.Range("K3").Value = "2016-09-03 01:07:30"

$wsZhCn.Cells.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0362aca070505afce837ff1ed03df425643cd9b7/e2e/$newFile1", "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8b1c13694558454dc05943e4acf3c9df0b6f6df4/e2e/$newFile1", "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0362aca070505afce837ff1ed03df425643cd9b7/e2e/$newFile2", "", "", $newFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8b1c13694558454dc05943e4acf3c9df0b6f6df4/e2e/$newFile2", "", "", $newFile2)

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("I3").Value = $newFile2

$deXlf = "224a5908-ba5e-498b-a323-6945e0c4307d.974b29d69231ec61ba9d9827bd1e821688ad796b.de-de.xlf"
$wsDeDe.Range("G2").Value = $deXlf
$wsDeDe.Range("J2").Value = $deXlf
$wsDeDe.Range("G3").Value = $deXlf
$wsDeDe.Range("J3").Value = $deXlf

$wsDeDe.Range("H2").Value = "2016-09-03 01:07:11"
$wsDeDe.Range("H3").Value = "2016-09-03 01:07:11"
$wsDeDe.Range("K2").Value = "2016-09-03 01:07:37"
$wsDeDe.Range("K3").Value = "2016-09-03 01:07:37"

$wsDeDe.Cells.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0362aca070505afce837ff1ed03df425643cd9b7/e2e/$newFile1", "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/04743fce537fdaa61a96cc516cdec0acb618a841/e2e/$newFile1", "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0362aca070505afce837ff1ed03df425643cd9b7/e2e/$newFile2", "", "", $newFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/04743fce537fdaa61a96cc516cdec0acb618a841/e2e/$newFile2", "", "", $newFile2)
